$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (64) of price data, following the same pattern used
# for the preceding rows (56-63): date stored as a text string, the
# remaining columns as plain numbers.
$ws.Cells.Item(64, 1).Value = "2024-10-03 00:00:00"
$ws.Cells.Item(64, 2).Value = 75650
$ws.Cells.Item(64, 3).Value = 10756.89
$ws.Cells.Item(64, 4).Value = 9519.370000000001
$ws.Cells.Item(64, 5).Value = 7.0457
